{"js": "// Apply the review-copy rewrite described by the diff:\n//  - Drop the \" | Review and Gameplay\" suffix from the title (both the\n//    H1 heading and the bold \"title\" line near the end of the doc).\n//  - Rewrite the \"What we like\" / \"What we don't like\" bullet list items.\n//  - Rewrite the italic meta-description line near the end of the doc.\n\nconst body = context.document.body;\n\n// Exact-text replacements. Each `search` + per-hit `insertText(\"Replace\")`\n// keeps the original run's formatting (bold/italic/heading style) intact,\n// since only the text inside the already-matched range is swapped out.\nconst replacements = [\n  [\n    \"Play East Sea Dragon King for Free | Review and Gameplay\",\n    \"Play East Sea Dragon King for Free\",\n  ],\n  [\n    \"Cluster payouts provide more chances to win\",\n    \"Cluster payouts for more winning opportunities\",\n  ],\n  [\n    \"High volatility adds excitement to gameplay\",\n    \"Expanding wilds and respins for dynamic gameplay\",\n  ],\n  [\n    \"Expanding wilds and respins feature increase chances of winning big\",\n    \"High volatility for thrilling gaming experience\",\n  ],\n  [\n    \"Beautifully designed symbols and graphics\",\n    \"Opportunity to win up to 6,000 times the total bet\",\n  ],\n  [\n    \"High volatility can be risky for inexperienced players\",\n    \"High volatility may not suit all players\",\n  ],\n  [\n    \"Not as many features as other NetEnt games\",\n    \"Limited bonus features\",\n  ],\n  [\n    \"Get a chance to win up to 6,000 times your bet with East Sea Dragon King. Play for free and learn about its features in our comprehensive review.\",\n    \"Read our review of East Sea Dragon King slot game and play for free. Discover its features and winning potential.\",\n  ],\n];\n\nconst searchResultsList = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true })\n);\nsearchResultsList.forEach((results) => results.load(\"items\"));\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const results = searchResultsList[i];\n  results.items.forEach((range) => range.insertText(newText, \"Replace\"));\n}\n\nawait context.sync();\n", "ps1": "# Apply the review-copy rewrite described by the diff:\n#  - Drop the \" | Review and Gameplay\" suffix from the title (both the\n#    H1 heading and the bold \"title\" line near the end of the doc).\n#  - Rewrite the \"What we like\" / \"What we don't like\" bullet list items.\n#  - Rewrite the italic meta-description line near the end of the doc.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText {\n    param(\n        [string]$OldText,\n        [string]$NewText\n    )\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $OldText\n    $find.Replacement.Text = $NewText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    # Replace every occurrence in the document (wdReplaceAll = 2) \u2014 the\n    # title string shows up twice (heading + bold line near the end) and\n    # both instances change to the same new text.\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nReplace-ExactText \"Play East Sea Dragon King for Free | Review and Gameplay\" \"Play East Sea Dragon King for Free\"\nReplace-ExactText \"Cluster payouts provide more chances to win\" \"Cluster payouts for more winning opportunities\"\nReplace-ExactText \"High volatility adds excitement to gameplay\" \"Expanding wilds and respins for dynamic gameplay\"\nReplace-ExactText \"Expanding wilds and respins feature increase chances of winning big\" \"High volatility for thrilling gaming experience\"\nReplace-ExactText \"Beautifully designed symbols and graphics\" \"Opportunity to win up to 6,000 times the total bet\"\nReplace-ExactText \"High volatility can be risky for inexperienced players\" \"High volatility may not suit all players\"\nReplace-ExactText \"Not as many features as other NetEnt games\" \"Limited bonus features\"\nReplace-ExactText \"Get a chance to win up to 6,000 times your bet with East Sea Dragon King. Play for free and learn about its features in our comprehensive review.\" \"Read our review of East Sea Dragon King slot game and play for free. Discover its features and winning potential.\"\n"}
